$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AZ1").Value = 0.91975406281535488
$ws.Range("BP1").Value = 0.69088698842695806
$ws.Range("C1").Value = 0.88818457506703719
$ws.Range("BJ2").Value = 0.71162190037889483
$ws.Range("BP2").Value = 0.81931950067113424
$ws.Range("D2").Value = 0.83770002496798712
$ws.Range("B3").Value = 0.76264433358092087
$ws.Range("D3").Value = 0.72330474189067939
$ws.Range("E3").Value = 0.97476964188683835
$ws.Range("I3").Value = 0.68632362437770955
$ws.Range("E4").Value = 0.94545971578703369
$ws.Range("F4").Value = 0.89834786704977931
$ws.Range("BL6").Value = 0.70118028540126209
$ws.Range("E6").Value = 0.79531997240073293
$ws.Range("BF7").Value = 0.96826157771932075
$ws.Range("E7").Value = 0.84322242827432614
$ws.Range("F7").Value = 0.86050829903361969
$ws.Range("H7").Value = 0.83173609470435572
$ws.Range("Q7").Value = 0.78107859466354057
$ws.Range("I8").Value = 0.92604456643097643
$ws.Range("X8").Value = 0.76288748067426293
$ws.Range("G9").Value = 0.75021505844618086
$ws.Range("H10").Value = 0.98787871966049412
$ws.Range("I10").Value = 0.9405812288715435
$ws.Range("K10").Value = 0.88414542969723586
$ws.Range("AR11").Value = 0.94949944475390624
$ws.Range("M11").Value = 0.73515887663867341
$ws.Range("J12").Value = 0.84807277193024388
$ws.Range("K12").Value = 0.88164477576356881
$ws.Range("M12").Value = 0.50694473924829375
$ws.Range("P13").Value = 0.85610344400150773
$ws.Range("M15").Value = 0.68622787893596637
$ws.Range("N15").Value = 0.91872003366055632
$ws.Range("P15").Value = 0.99007883211269565
$ws.Range("S15").Value = 0.60013942726850911
$ws.Range("AI16").Value = 0.84573297181557827
$ws.Range("BC16").Value = 0.64897710230782546
$ws.Range("N16").Value = 0.84546077644848627
$ws.Range("Q16").Value = 0.83051995434604375
$ws.Range("O17").Value = 0.80264925778091545
$ws.Range("S17").Value = 0.82230567950177802
$ws.Range("P18").Value = 0.74858739419265152
$ws.Range("S18").Value = 0.69048197081921836
$ws.Range("R20").Value = 0.52900046657056421
$ws.Range("S20").Value = 0.86468392682999673
$ws.Range("U20").Value = 0.97791007565900956
$ws.Range("AO22").Value = 0.50906738229022075
$ws.Range("U22").Value = 0.95770660864229207
$ws.Range("AM23").Value = 0.95695871465503135
$ws.Range("L23").Value = 0.93003540703551835
$ws.Range("U23").Value = 0.77505958765227789
$ws.Range("V24").Value = 0.7779437197230048
$ws.Range("W24").Value = 0.86700214020334454
$ws.Range("AA25").Value = 0.84161567392931114
$ws.Range("Z25").Value = 0.94920790480382888
$ws.Range("AA26").Value = 0.87816663426183594
$ws.Range("AB27").Value = 0.54528711969093158
$ws.Range("AR27").Value = 0.80469882820150107
$ws.Range("Z28").Value = 0.86943549206766746
$ws.Range("AB29").Value = 0.89646914390960819
$ws.Range("AD29").Value = 0.87653925829346413
$ws.Range("AB30").Value = 0.86436664106678085
$ws.Range("D30").Value = 0.69813828714467063
$ws.Range("AC31").Value = 0.82642874261569665
$ws.Range("AU31").Value = 0.67201957450397953
$ws.Range("AE32").Value = 0.70895050737275067
$ws.Range("AF33").Value = 0.59317674241057428
$ws.Range("AH33").Value = 0.88347439474527745
$ws.Range("AI33").Value = 0.8132977225522664
$ws.Range("AF34").Value = 0.94307151885699225
$ws.Range("AH35").Value = 0.65325712622845877
$ws.Range("AK36").Value = 0.84854541978349274
$ws.Range("AI37").Value = 0.75579254283358233
$ws.Range("AM37").Value = 0.90155018669134013
$ws.Range("M37").Value = 0.82360405313901608
$ws.Range("AM38").Value = 0.82141349094812766
$ws.Range("AD39").Value = 0.82871316105428194
$ws.Range("AN39").Value = 0.90780976577027639
$ws.Range("AO39").Value = 0.63099555163171916
$ws.Range("T39").Value = 0.93943902327662332
$ws.Range("AL40").Value = 0.83716787919587843
$ws.Range("AN42").Value = 0.87730004247641058
$ws.Range("AO42").Value = 0.81612142809158872
$ws.Range("AQ42").Value = 0.77473998540627209
$ws.Range("AO43").Value = 0.81303906911178836
$ws.Range("AP44").Value = 0.80322116124231702
$ws.Range("AQ44").Value = 0.97636113927199752
$ws.Range("AQ45").Value = 0.68528960678509954
$ws.Range("AU45").Value = 0.98927120892229192
$ws.Range("AR46").Value = 0.80117816949852605
$ws.Range("AS46").Value = 0.90649516077659853
$ws.Range("AU46").Value = 0.60110715767763523
$ws.Range("BK47").Value = 0.73005529751257159
$ws.Range("AT48").Value = 0.72642366769890643
$ws.Range("AW48").Value = 0.73582615060134937
$ws.Range("AN49").Value = 0.67240766380474981
$ws.Range("AU49").Value = 0.72377708433468324
$ws.Range("AV50").Value = 0.92955951659171798
$ws.Range("AW50").Value = 0.92428756118193911
$ws.Range("AY50").Value = 0.85730485250500688
$ws.Range("AZ50").Value = 0.82202774071993345
$ws.Range("AW51").Value = 0.87807102674456794
$ws.Range("BA51").Value = 0.72646540882458954
$ws.Range("BH51").Value = 0.7183226716180735
$ws.Range("BO51").Value = 0.90626575079154503
$ws.Range("AH52").Value = 0.89961531110530035
$ws.Range("AZ53").Value = 0.97974168521642291
$ws.Range("BC53").Value = 0.88246764810152611
$ws.Range("BJ53").Value = 0.84977854657006979
$ws.Range("AC54").Value = 0.98466718804296804
$ws.Range("BA54").Value = 0.96534632256426123
$ws.Range("BM54").Value = 0.93186174303007885
$ws.Range("BE56").Value = 0.99898599692008538
$ws.Range("BF56").Value = 0.82390067560935787
$ws.Range("BC57").Value = 0.73580808263786202
$ws.Range("BG57").Value = 0.99704757672321076
$ws.Range("BN57").Value = 0.78624459981834627
$ws.Range("AZ58").Value = 0.83906215185973188
$ws.Range("BE58").Value = 0.9449768482713321
$ws.Range("BH58").Value = 0.89893436904677193
$ws.Range("BH59").Value = 0.86718029319361423
$ws.Range("BI59").Value = 0.69641256659084483
$ws.Range("BI60").Value = 0.9959562308782175
$ws.Range("X60").Value = 0.99662288591435966
$ws.Range("AE61").Value = 0.7760779000536957
$ws.Range("BK61").Value = 0.95735672772522373
$ws.Range("AJ64").Value = 0.89521842884066738
$ws.Range("BK64").Value = 0.75004755833844583
$ws.Range("BK65").Value = 0.98118341288394118
$ws.Range("BL65").Value = 0.7869614567748755
$ws.Range("BL66").Value = 0.72951265076962379
$ws.Range("BM66").Value = 0.79445399380035275
$ws.Range("A67").Value = 0.88065064194754927
$ws.Range("BM67").Value = 0.99991956515623914
$ws.Range("BN68").Value = 0.75726905094319585
$ws.Range("BO68").Value = 0.84502281112698907
